$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.99999999005685725
$ws.Range("A2").Value = 0.99530002778030868
$ws.Range("A3").Value = 0.97387034383030491
$ws.Range("A4").Value = 0.96640767468129507
$ws.Range("A5").Value = 0.95943385547291871
$ws.Range("A6").Value = 0.9452216443019148
$ws.Range("A7").Value = 0.94375267116972061
$ws.Range("A8").Value = 0.94456077982082287
$ws.Range("A9").Value = 0.94891496341512571
$ws.Range("A10").Value = 0.95430610660932613
$ws.Range("A11").Value = 0.95373138154321824
$ws.Range("A12").Value = 0.95095240398346492
$ws.Range("A13").Value = 0.93966816109022577
$ws.Range("A14").Value = 0.9355014956995954
$ws.Range("A15").Value = 0.93291036898897839
$ws.Range("A16").Value = 0.9304041168518038
$ws.Range("A17").Value = 0.92669645780111143
$ws.Range("A18").Value = 0.9255875773683786
$ws.Range("A19").Value = 0.99424687785984456
$ws.Range("A20").Value = 0.98713004410619454
$ws.Range("A21").Value = 0.98573159426207035
$ws.Range("A22").Value = 0.98446709850918479
$ws.Range("A23").Value = 0.96683487342347507
$ws.Range("A24").Value = 0.95381329791327696
$ws.Range("A25").Value = 0.94735619893977008
$ws.Range("A26").Value = 0.94235730747869662
$ws.Range("A27").Value = 0.93911517276022582
$ws.Range("A28").Value = 0.92530038516936652
$ws.Range("A29").Value = 0.91570934547008287
$ws.Range("A30").Value = 0.91188614568656379
$ws.Range("A31").Value = 0.91150334280266365
$ws.Range("A32").Value = 0.91237438819348304
$ws.Range("A33").Value = 0.91185439506984234
